$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing hours value for the existing last row (row 29)
$ws.Range("C29").Value = 6

# New row 30
$ws.Range("A30").Value = 45663
$ws.Range("A30").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("B30").Value = "ART"
$ws.Range("C30").Value = 6

# New row 31
$ws.Range("A31").Value = 45664
$ws.Range("A31").NumberFormat = $ws.Range("A29").NumberFormat
$ws.Range("B31").Value = "more writing and planning - reworked one of the characters and the beginning of the story"
$ws.Range("C31").Value = 5

# Update the view to match the committed state
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("E30").Select()
